$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (this shifts rows 5-19 down to 6-20,
# and carries each row's existing formatting/styles down with it).
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new data point (x = 0.01)
# Match the formatting of the row directly below it (old row 5, now row 6)
# by copying its style, then overwrite with the numeric values.
$ws.Range("C6:G6").Copy()
$ws.Range("C5:G5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C5").Value = 0.01
$ws.Range("D5").Value = 0.9872
$ws.Range("E5").Value = 0.9861
$ws.Range("F5").Value = 0.9882
$ws.Range("G5").Value = 0.9871

# Update the selection to match the post-edit state recorded in the file.
$ws.Range("I12").Select()
